$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
# "Volume 30   Number  48" -> "Volume 30   Number  49"
$ws.Range("A8").Value = "Volume 30   Number  49"
# "Report Covering the Week  11/27/2023  Through  12/3/2023"
#   -> "Report Covering the Week  12/4/2023  Through  12/10/2023"
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# --- Row 20: D20/E20 switch from blank-marker text to real numbers ---
# Reuse the existing number formats already used by sibling cells so the
# written cells land on the same style records as the target file.
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 3
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E20").Value = -100

$ws.Range("N14").Value = -94.117647058823
$ws.Range("J15").Value = 19
$ws.Range("K15").Value = -5.263157894736
$ws.Range("M15").Value = 38.461538461538
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 80
$ws.Range("I16").Value = 181
$ws.Range("J16").Value = 212
$ws.Range("K16").Value = -14.622641509434
$ws.Range("L16").Value = 1.117318435754
$ws.Range("M16").Value = 2.259887005649
$ws.Range("N16").Value = -80.641711229946
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 58.333333333333
$ws.Range("I17").Value = 212
$ws.Range("J17").Value = 195
$ws.Range("K17").Value = 8.717948717948
$ws.Range("L17").Value = 10.416666666666
$ws.Range("M17").Value = 46.206896551724
$ws.Range("N17").Value = -61.732851985559
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 42.105263157894
$ws.Range("I18").Value = 265
$ws.Range("J18").Value = 360
$ws.Range("K18").Value = -26.388888888888
$ws.Range("L18").Value = -4.676258992805
$ws.Range("M18").Value = 26.794258373205
$ws.Range("N18").Value = -70.090293453724
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = 22.222222222222
$ws.Range("F19").Value = 67
$ws.Range("G19").Value = 77
$ws.Range("H19").Value = -12.987012987013
$ws.Range("I19").Value = 929
$ws.Range("J19").Value = 985
$ws.Range("K19").Value = -5.685279187817
$ws.Range("L19").Value = 21.121251629726
$ws.Range("M19").Value = 26.912568306010
$ws.Range("N19").Value = -41.901188242651
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -60
$ws.Range("J20").Value = 48
$ws.Range("K20").Value = -12.5
$ws.Range("M20").Value = -10.638297872340
$ws.Range("N20").Value = -92.030360531309
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 13.793103448275
$ws.Range("F21").Value = 133
$ws.Range("G21").Value = 125
$ws.Range("H21").Value = 6.4
$ws.Range("I21").Value = 1648
$ws.Range("J21").Value = 1825
$ws.Range("K21").Value = -9.698630136986
$ws.Range("L21").Value = 7.712418300653
$ws.Range("M21").Value = 24.189902034664
$ws.Range("N21").Value = -63.804085218537
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -37.5
$ws.Range("J23").Value = 127
$ws.Range("K23").Value = -9.448818897637
$ws.Range("L23").Value = -40.414507772020
$ws.Range("M23").Value = 6.481481481481
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = -1.612903225806
$ws.Range("I24").Value = 1424
$ws.Range("J24").Value = 2069
$ws.Range("K24").Value = -31.174480425326
$ws.Range("L24").Value = 24.475524475524
$ws.Range("M24").Value = -13.801452784503
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 140
$ws.Range("F25").Value = 47
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = 147.368421052632
$ws.Range("I25").Value = 451
$ws.Range("J25").Value = 436
$ws.Range("K25").Value = 3.440366972477
$ws.Range("L25").Value = 18.997361477572
$ws.Range("M25").Value = 1.121076233183
$ws.Range("H26").Value = -100
$ws.Range("J26").Value = 34
$ws.Range("K26").Value = -5.882352941176
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 40
$ws.Range("I27").Value = 62
$ws.Range("K27").Value = -28.735632183908
$ws.Range("L27").Value = -6.060606060606
$ws.Range("L28").Value = -50
$ws.Range("L29").Value = -40

# --- Row 26: F26 switches from a number to the blank-marker text "0" ---
# Force text entry via a Text-formatted cell, then pull the exact
# format (and style record) from C20, which already displays as "0".
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "0"
$ws.Range("C20").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
